$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.988.25'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.24'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5186'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.88%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2820'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +4.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.66'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.25%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.754.81'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07026'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.523'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '77.41'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.996.22'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.50'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006612'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.978.92'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.151'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.648'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.145'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '139.65'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.510'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.10'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '102.44'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08277'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.668'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.441'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04481'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.608'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9888'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6171'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.676'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01586'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.926'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.13%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '100.24'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.77%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7396'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.070'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05457'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.349'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1128'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.08'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '30.13'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.624'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.53%  '
